$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column).
# This shifts the old N,O,P data into O,P,Q and copies the
# formatting (width/style) of the column to the left (M).
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Select a cell on the repayment schedule sheet (mirrors the author's
# on-screen selection after editing) and make it the active sheet/tab.
# Activating this sheet automatically drops tabSelected from whichever
# sheet previously held it (the last tab, "Acc_Periodic").
$ws.Activate()
$ws.Range("S6").Select()
